# Corrección a Diebold Mariano y revisión de Cap1
# Update DM_Stat (col C) and P_Value (col D) for rows 2-11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.3325567634943948
$ws.Range("D2").Value = 0.7426168946456271

$ws.Range("C3").Value = -0.2488684190128471
$ws.Range("D3").Value = 0.8057711849165643

$ws.Range("C4").Value = 0.3499099854345538
$ws.Range("D4").Value = 0.729733879251119

$ws.Range("C5").Value = 0.02638432778644234
$ws.Range("D5").Value = 0.9791886644828323

$ws.Range("C6").Value = 0.163690994870002
$ws.Range("D6").Value = 0.8714688851860588

$ws.Range("C7").Value = 0.8517170263947755
$ws.Range("D7").Value = 0.4035475488691964

$ws.Range("C8").Value = 0.5157727461775639
$ws.Range("D8").Value = 0.6111579654143564

$ws.Range("C9").Value = 0.7384218225459965
$ws.Range("D9").Value = 0.4680630473385654

$ws.Range("C10").Value = 0.2497552347405249
$ws.Range("D10").Value = 0.8050939345185331

$ws.Range("C11").Value = -0.2761774934693786
$ws.Range("D11").Value = 0.7849895650422689
